# update scripts wuth new tpm
# NATMI re-run with refreshed TPM input for the "ECs" cluster changed the
# ligand (Tnfsf13) and receptor (Fas) average/total expression values for
# that cluster, which ripple into every derived specificity / edge-weight
# column (I, J, O, P, Q, R, S, T) for every sending/target cluster pair.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.5347833333333334
$ws.Range("H2").Value = 1.60435
$ws.Range("I2").Value = 0.196822066153855
$ws.Range("J2").Value = 0.196822066153855
$ws.Range("M2").Value = 17.125047
$ws.Range("N2").Value = 51.375141
$ws.Range("O2").Value = 0.6137494773138626
$ws.Range("P2").Value = 0.6137494773138626
$ws.Range("Q2").Value = 9.15818971815
$ws.Range("R2").Value = 82.42370746335001
$ws.Range("S2").Value = 0.120799440225763
$ws.Range("T2").Value = 0.120799440225763

$ws.Range("G3").Value = 0.5347833333333334
$ws.Range("H3").Value = 1.60435
$ws.Range("I3").Value = 0.196822066153855
$ws.Range("J3").Value = 0.196822066153855
$ws.Range("O3").Value = 0.2391747492820853
$ws.Range("P3").Value = 0.2391747492820853
$ws.Range("Q3").Value = 3.568895470677778
$ws.Range("R3").Value = 32.1200592361
$ws.Range("S3").Value = 0.0470748683255303
$ws.Range("T3").Value = 0.04707486832553029

$ws.Range("G4").Value = 0.5347833333333334
$ws.Range("H4").Value = 1.60435
$ws.Range("I4").Value = 0.196822066153855
$ws.Range("J4").Value = 0.196822066153855
$ws.Range("O4").Value = 0.147075773404052
$ws.Range("P4").Value = 0.147075773404052
$ws.Range("Q4").Value = 2.194621560694444
$ws.Range("R4").Value = 19.75159404625
$ws.Range("S4").Value = 0.02894775760256172
$ws.Range("T4").Value = 0.02894775760256172

$ws.Range("I5").Value = 0.1891972429821067
$ws.Range("J5").Value = 0.1891972429821067
$ws.Range("M5").Value = 17.125047
$ws.Range("N5").Value = 51.375141
$ws.Range("O5").Value = 0.6137494773138626
$ws.Range("P5").Value = 0.6137494773138626
$ws.Range("Q5").Value = 8.803404411102
$ws.Range("R5").Value = 79.230639699918
$ws.Range("S5").Value = 0.1161197089894919
$ws.Range("T5").Value = 0.1161197089894919

$ws.Range("I6").Value = 0.1891972429821067
$ws.Range("J6").Value = 0.1891972429821067
$ws.Range("O6").Value = 0.2391747492820853
$ws.Range("P6").Value = 0.2391747492820853
$ws.Range("S6").Value = 0.04525120315510715
$ws.Range("T6").Value = 0.04525120315510715

$ws.Range("I7").Value = 0.1891972429821067
$ws.Range("J7").Value = 0.1891972429821067
$ws.Range("O7").Value = 0.147075773404052
$ws.Range("P7").Value = 0.147075773404052
$ws.Range("S7").Value = 0.0278263308375077
$ws.Range("T7").Value = 0.0278263308375077

$ws.Range("I8").Value = 0.6139806908640383
$ws.Range("J8").Value = 0.6139806908640382
$ws.Range("M8").Value = 17.125047
$ws.Range("N8").Value = 51.375141
$ws.Range("O8").Value = 0.6137494773138626
$ws.Range("P8").Value = 0.6137494773138626
$ws.Range("Q8").Value = 28.568705532327
$ws.Range("R8").Value = 257.118349790943
$ws.Range("S8").Value = 0.3768303280986078
$ws.Range("T8").Value = 0.3768303280986077

$ws.Range("I9").Value = 0.6139806908640383
$ws.Range("J9").Value = 0.6139806908640382
$ws.Range("O9").Value = 0.2391747492820853
$ws.Range("P9").Value = 0.2391747492820853
$ws.Range("S9").Value = 0.1468486778014479
$ws.Range("T9").Value = 0.1468486778014479

$ws.Range("I10").Value = 0.6139806908640383
$ws.Range("J10").Value = 0.6139806908640382
$ws.Range("O10").Value = 0.147075773404052
$ws.Range("P10").Value = 0.147075773404052
$ws.Range("S10").Value = 0.09030168496398261
$ws.Range("T10").Value = 0.0903016849639826
